# Adds a new forecast-date column (Z) for 2020-05-19 to both the "cases"
# and "deaths" tables, plus the corresponding new target-date row (38),
# and fills in the previously-missing "observed" value for row 24
# (target date 2020-05-05) now that it has become available.

$wb = $excel.ActiveWorkbook

# New forecast/target date introduced by this edit.
$newDate = "2020-05-19"

# Per-sheet data: B24 "observed" backfill value, and the new Z column
# values keyed by row number (rows 25-37 correspond to target dates
# 2020-05-06 .. 2020-05-18), plus the brand-new row 38 (target date
# 2020-05-19, same as the new column's forecast date).
$sheetData = @{
    "cases" = @{
        B24 = 34053
        Z = @{
            25 = 36076; 26 = 37781; 27 = 39316; 28 = 40595; 29 = 41528;
            30 = 42662; 31 = 43650; 32 = 44395; 33 = 45222; 34 = 45841;
            35 = 46729; 36 = 47237; 37 = 47859
        }
        Z38 = 48455
    }
    "deaths" = @{
        B24 = 2851
        Z = @{
            25 = 2990; 26 = 3110; 27 = 3218; 28 = 3304; 29 = 3360;
            30 = 3444; 31 = 3515; 32 = 3565; 33 = 3627; 34 = 3668;
            35 = 3739; 36 = 3771; 37 = 3817
        }
        Z38 = 3863
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $sheetData[$sheetName]

    # --- Header: Z1 gets the label for 2020-05-05, exactly like A24 ---
    # (copy A24 instead of typing the string so the date-like text isn't
    # auto-converted into a date value/format)
    $ws.Range("A24").Copy()
    $ws.Range("Z1").PasteSpecial(-4163)  # xlPasteValues

    # --- Backfill the now-known "observed" figure for row 24 ---
    $ws.Range("B24").Value = $data["B24"]

    # --- New forecast column Z for existing rows 25-37 ---
    foreach ($r in $data["Z"].Keys) {
        $ws.Range("Z$r").Value = $data["Z"][$r]
    }

    # --- New row 38 for target date 2020-05-19 ---
    # Build the new shared string via a throw-away formula cell (a
    # formula result is never re-interpreted as a date), copy it in as a
    # plain value, then clean up the helper cell.
    $ws.Range("AA1").Formula = '="' + $newDate + '"'
    $ws.Range("AA1").Copy()
    $ws.Range("A38").PasteSpecial(-4163)  # xlPasteValues
    $ws.Range("AA1").ClearContents()

    $ws.Range("Z38").Value = $data["Z38"]
}
